# Update gh-pages to output generated at 456a3b4
# Updates numeric counter values (column F) on the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F9").Value = 10658
$ws1.Range("F15").Value = 7524
$ws1.Range("F18").Value = 256

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F12").Value = 10658
$ws4.Range("F18").Value = 7524
$ws4.Range("F21").Value = 256
